$wb = $excel.ActiveWorkbook

# --- About sheet: update the "last updated" date (C1) ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45379

# --- RAF-capacity sheet: bump the hydrogen RAF multipliers to 1 ---
$wsCapacity = $wb.Worksheets.Item("RAF-capacity")
$wsCapacity.Range("B24").Value = 1
$wsCapacity.Range("B25").Value = 1

# Widen column A slightly so the hydrogen technology labels fit
$wsCapacity.Columns.Item(1).ColumnWidth = 28.14

# Make RAF-capacity the active/selected sheet, zoomed in, with B25 selected
$wsCapacity.Activate()
$wsCapacity.Range("B25").Select()
$excel.ActiveWindow.Zoom = 80
